$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7907.6
$ws.Range("I19").Value = 955.8570999999999
$ws.Range("J19").Value = 11650.846
$ws.Range("K19").Value = 955.8570999999999
$ws.Range("L19").Value = 11650.846
$ws.Range("M19").Value = -780.8570999999999
$ws.Range("N19").Value = -12000.846
$ws.Range("H29").Value = 3483.2856
$ws.Range("I29").Value = 1794.3334
$ws.Range("J29").Value = 4750
$ws.Range("K29").Value = 5383.0002
$ws.Range("L29").Value = 14250
$ws.Range("M29").Value = -5102.0002
$ws.Range("N29").Value = -14812
$ws.Range("H40").Value = 4882.1665
$ws.Range("I40").Value = 5258.8
$ws.Range("J40").Value = 2999
$ws.Range("K40").Value = 5258.8
$ws.Range("L40").Value = 2999
$ws.Range("M40").Value = -5083.8
$ws.Range("N40").Value = -3349
$ws.Range("H88").Value = 13893858
$ws.Range("I88").Value = 41670604
$ws.Range("J88").Value = 5484
$ws.Range("K88").Value = 41670604
$ws.Range("L88").Value = 5484
$ws.Range("M88").Value = -41670198
$ws.Range("N88").Value = -6296
$ws.Range("H91").Value = 13893858
$ws.Range("I91").Value = 41670604
$ws.Range("J91").Value = 5484
$ws.Range("K91").Value = 41670604
$ws.Range("L91").Value = 5484
$ws.Range("M91").Value = -41669200
$ws.Range("N91").Value = -8292
$ws.Range("H116").Value = 3938.7778
$ws.Range("I116").Value = 3739.8
$ws.Range("J116").Value = 4015.3076
$ws.Range("K116").Value = 3739.8
$ws.Range("L116").Value = 4015.3076
$ws.Range("M116").Value = -297.8000000000002
$ws.Range("N116").Value = -10899.3076

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5404.7334
$ws.Range("I61").Value = 2869.875
$ws.Range("K61").Value = 2869.875
$ws.Range("M61").Value = -2657.875
$ws.Range("H88").Value = 1835.091
$ws.Range("J88").Value = 1854.9474
$ws.Range("L88").Value = 1854.9474
$ws.Range("N88").Value = -2666.9474
$ws.Range("H91").Value = 1835.091
$ws.Range("J91").Value = 1854.9474
$ws.Range("L91").Value = 1854.9474
$ws.Range("M91").Value = -666.9474
$ws.Range("N91").Value = -4662.9474
$ws.Range("H97").Value = 8841.362999999999
$ws.Range("I97").Value = 9756.875
$ws.Range("J97").Value = 6400
$ws.Range("K97").Value = 9756.875
$ws.Range("L97").Value = 6400
$ws.Range("M97").Value = -9260.875
$ws.Range("N97").Value = -7392
$ws.Range("H136").Value = 5404.7334
$ws.Range("I136").Value = 2869.875
$ws.Range("K136").Value = 8609.625
$ws.Range("M136").Value = -6059.625

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 6174.5
$ws.Range("I97").Value = 6174.5
$ws.Range("K97").Value = 6174.5
$ws.Range("M97").Value = -5183.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7733
$ws.Range("I16").Value = 7224.75
$ws.Range("J16").Value = 8749.5
$ws.Range("K16").Value = 7224.75
$ws.Range("L16").Value = 8749.5
$ws.Range("M16").Value = -6937.75
$ws.Range("N16").Value = -9323.5
$ws.Range("H58").Value = 6167.35
$ws.Range("I58").Value = 6112.84
$ws.Range("J58").Value = 6258.2
$ws.Range("K58").Value = 6112.84
$ws.Range("L58").Value = 6258.2
$ws.Range("M58").Value = -5909.84
$ws.Range("N58").Value = -6664.2
$ws.Range("H94").Value = 5439.467
$ws.Range("J94").Value = 4395.8
$ws.Range("L94").Value = 4395.8
$ws.Range("N94").Value = -5297.8
$ws.Range("H113").Value = 7733
$ws.Range("I113").Value = 7224.75
$ws.Range("J113").Value = 8749.5
$ws.Range("K113").Value = 7224.75
$ws.Range("L113").Value = 8749.5
$ws.Range("M113").Value = -5054.75
$ws.Range("N113").Value = -13089.5
$ws.Range("H132").Value = 6618.484
$ws.Range("I132").Value = 7833.913
$ws.Range("J132").Value = 3124.125
$ws.Range("K132").Value = 23501.739
$ws.Range("L132").Value = 9372.375
$ws.Range("M132").Value = -20971.739
$ws.Range("N132").Value = -14432.375
$ws.Range("H136").Value = 6167.35
$ws.Range("I136").Value = 6112.84
$ws.Range("J136").Value = 6258.2
$ws.Range("K136").Value = 18338.52
$ws.Range("L136").Value = 18774.6
$ws.Range("M136").Value = -15788.52
$ws.Range("N136").Value = -23874.6

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25508.5
$ws.Range("J57").Value = 29999.2
$ws.Range("L57").Value = 29999.2
$ws.Range("N57").Value = -31639.2
$ws.Range("H70").Value = 3431.0688
$ws.Range("I70").Value = 3179.0588
$ws.Range("K70").Value = 3179.0588
$ws.Range("M70").Value = -2909.0588
$ws.Range("H73").Value = 3431.0688
$ws.Range("I73").Value = 3179.0588
$ws.Range("K73").Value = 3179.0588
$ws.Range("M73").Value = -2243.0588
$ws.Range("H102").Value = 5012.8945
$ws.Range("I102").Value = 5957
$ws.Range("J102").Value = 2967.3333
$ws.Range("K102").Value = 5957
$ws.Range("L102").Value = 2967.3333
$ws.Range("M102").Value = -4335
$ws.Range("N102").Value = -6211.3333
$ws.Range("H113").Value = 11297.333
$ws.Range("I113").Value = 3957.2
$ws.Range("K113").Value = 3957.2
$ws.Range("M113").Value = -1787.2
$ws.Range("H126").Value = 6046.9443
$ws.Range("I126").Value = 5913.0835
$ws.Range("J126").Value = 6314.6665
$ws.Range("K126").Value = 17739.2505
$ws.Range("L126").Value = 18943.9995
$ws.Range("M126").Value = -15269.2505
$ws.Range("N126").Value = -23883.9995

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3901.5715
$ws.Range("I7").Value = 3353.875
$ws.Range("J7").Value = 4631.8335
$ws.Range("K7").Value = 3353.875
$ws.Range("L7").Value = 4631.8335
$ws.Range("M7").Value = -3241.875
$ws.Range("N7").Value = -4855.8335
$ws.Range("H46").Value = 2310.8572
$ws.Range("I46").Value = 761.1111
$ws.Range("J46").Value = 2847.3076
$ws.Range("K46").Value = 761.1111
$ws.Range("L46").Value = 2847.3076
$ws.Range("M46").Value = -573.1111
$ws.Range("N46").Value = -3223.3076
$ws.Range("H126").Value = 3901.5715
$ws.Range("I126").Value = 3353.875
$ws.Range("J126").Value = 4631.8335
$ws.Range("K126").Value = 10061.625
$ws.Range("L126").Value = 13895.5005
$ws.Range("M126").Value = -7591.625
$ws.Range("N126").Value = -18835.5005
$ws.Range("H132").Value = 42901.793
$ws.Range("I132").Value = 57126.668
$ws.Range("J132").Value = 5561.5
$ws.Range("K132").Value = 171380.004
$ws.Range("L132").Value = 16684.5
$ws.Range("M132").Value = -168850.004
$ws.Range("N132").Value = -21744.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1554.2
$ws.Range("I96").Value = 1270.6
$ws.Range("K96").Value = 1270.6
$ws.Range("M96").Value = 102.4000000000001
$ws.Range("H100").Value = 915.7
$ws.Range("I100").Value = 766.7143
$ws.Range("J100").Value = 1263.3334
$ws.Range("K100").Value = 1533.4286
$ws.Range("L100").Value = 2526.6668
$ws.Range("M100").Value = -992.4286
$ws.Range("N100").Value = -3608.6668
$ws.Range("H113").Value = 697.6
$ws.Range("I113").Value = 668.4
$ws.Range("J113").Value = 726.8
$ws.Range("K113").Value = 2005.2
$ws.Range("L113").Value = 2180.4
$ws.Range("M113").Value = 164.8000000000002
$ws.Range("N113").Value = -6520.4
$ws.Range("H132").Value = 2119.7144
$ws.Range("I132").Value = 1761
$ws.Range("K132").Value = 5283
$ws.Range("M132").Value = -2753
$ws.Range("H136").Value = 20570.857
$ws.Range("J136").Value = 7798.4
$ws.Range("L136").Value = 23395.2
$ws.Range("N136").Value = -28495.2
